# Appends new newsbot rows (22-31) to the Historico sheet, covering columns A-H.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22
$ws.Cells.Item(22, 1).Value = "04/01/2026 22:44:52"
$ws.Cells.Item(22, 2).Value = "04/01 18:14"
$ws.Cells.Item(22, 3).Value = "BBC Brasil"
$ws.Cells.Item(22, 4).Value = "Trump diz que sucessora de Maduro 'pagará preço muito alto' se 'não fizer o certo' na Venezuela; o que aconteceu até agora após ataque americano"
$ws.Cells.Item(22, 5).Value = "https://www.bbc.com/portuguese/articles/cm2403jvm03o?at_medium=RSS&at_campaign=rss"
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = "sc"
$ws.Cells.Item(22, 8).Value = "Ação militar e prisão de Maduro marcam auge de meses de e<b>sc</b>alada de tensão entre os dois países. Maioria dos líderes da América Latina cond"

# Row 23
$ws.Cells.Item(23, 1).Value = "04/01/2026 22:44:53"
$ws.Cells.Item(23, 2).Value = "04/01 18:00"
$ws.Cells.Item(23, 3).Value = "Folha de S.Paulo - Poder - Principal"
$ws.Cells.Item(23, 4).Value = "Republicanos se distancia de defesa da pena de morte feita por ministro"
$ws.Cells.Item(23, 5).Value = "https://redir.folha.com.br/redir/online/poder/rss091/*https://www1.folha.uol.com.br/colunas/painel/2026/01/republicanos-se-distancia-de-defesa-da-pena-de-morte-feita-por-ministro.shtml"
$ws.Cells.Item(23, 6).Value = 2
$ws.Cells.Item(23, 7).Value = "lula"
$ws.Cells.Item(23, 8).Value = "orte para autores de feminicídio constrangeu integrantes do governo Luiz Inácio Lula da Silva e surpreendeu até integrantes de seu partido, o Republicanos. Segundo "

# Row 24
$ws.Cells.Item(24, 1).Value = "04/01/2026 22:44:54"
$ws.Cells.Item(24, 2).Value = "04/01 18:00"
$ws.Cells.Item(24, 3).Value = "Folha de S.Paulo - Mercado - Principal"
$ws.Cells.Item(24, 4).Value = "O que realmente merece planejamento financeiro"
$ws.Cells.Item(24, 5).Value = "https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/blogs/de-grao-em-grao/2026/01/o-que-realmente-merece-planejamento-financeiro.shtml"
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = "sc"
$ws.Cells.Item(24, 8).Value = "nsformar a vida em uma planilha cansativa. O efeito dessa confusão é curioso: di<b>sc</b>ute-se demais o que é pequeno e decide-se de menos o que é grande.`n&lt;a href=&quot;http"

# Row 25
$ws.Cells.Item(25, 1).Value = "04/01/2026 22:44:55"
$ws.Cells.Item(25, 2).Value = "04/01 15:00"
$ws.Cells.Item(25, 3).Value = "Folha de S.Paulo - Poder - Principal"
$ws.Cells.Item(25, 4).Value = "Pacheco trata indicação ao Supremo como página virada, dizem aliados"
$ws.Cells.Item(25, 5).Value = "https://redir.folha.com.br/redir/online/poder/rss091/*https://www1.folha.uol.com.br/colunas/painel/2026/01/pacheco-trata-indicacao-ao-supremo-como-pagina-virada-dizem-aliados.shtml"
$ws.Cells.Item(25, 6).Value = 4
$ws.Cells.Item(25, 7).Value = "senado"
$ws.Cells.Item(25, 8).Value = "O ex-presidente do &lt;a href=&quot;https://www1.folha.uol.com.br/folha-topicos/<b>senado</b>/&quot;&gt;Senado&lt;/a&gt; &lt;a href=&quot;https://www1.folha.uol.com.br/folha-topicos/rodrigo-pache"

# Row 26
$ws.Cells.Item(26, 1).Value = "04/01/2026 22:44:56"
$ws.Cells.Item(26, 2).Value = "04/01 15:00"
$ws.Cells.Item(26, 3).Value = "Folha de S.Paulo - Mercado - Principal"
$ws.Cells.Item(26, 4).Value = "Portal da reforma tributária permite monitorar economia em tempo real"
$ws.Cells.Item(26, 5).Value = "https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/blogs/que-imposto-e-esse/2026/01/portal-da-reforma-tributaria-permite-monitorar-economia-em-tempo-real.shtml"
$ws.Cells.Item(26, 6).Value = 2
$ws.Cells.Item(26, 7).Value = "imposto"
$ws.Cells.Item(26, 8).Value = "lha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/blogs/que-<b>imposto</b>-e-esse/2026/01/portal-da-reforma-tributaria-permite-monitorar-economia-em-tempo"

# Row 27
$ws.Cells.Item(27, 1).Value = "04/01/2026 22:45:00"
$ws.Cells.Item(27, 3).Value = "VEJA"
$ws.Cells.Item(27, 4).Value = "Itamaraty reforça na Celac preocupação com captura de Maduro na Venezuela"
$ws.Cells.Item(27, 5).Value = "https://veja.abril.com.br/mundo/itamaraty-reforca-na-celac-preocupacao-com-captura-de-maduro-na-venezuela/"
$ws.Cells.Item(27, 6).Value = 2
$ws.Cells.Item(27, 7).Value = "lula"
$ws.Cells.Item(27, 8).Value = "tro das Relações Exteriores, Mauro Vieira reiterou posicionamento do presidente Lula em reunião com países latino-americanos e caribenhos"

# Row 28
$ws.Cells.Item(28, 1).Value = "04/01/2026 22:45:01"
$ws.Cells.Item(28, 3).Value = "VEJA"
$ws.Cells.Item(28, 4).Value = "Forças armadas da Venezuela reconhecem Delcy Rodríguez como presidente interina"
$ws.Cells.Item(28, 5).Value = "https://veja.abril.com.br/mundo/forcas-armadas-da-venezuela-reconhecem-delcy-rodriguez-como-presidente-interina/"
$ws.Cells.Item(28, 6).Value = 2
$ws.Cells.Item(28, 7).Value = "câmara"
$ws.Cells.Item(28, 8).Value = "O chefe do exército venezuelano referendou a decisão da Câmara Constitucional da Suprema Corte da Venezuela"

# Row 29
$ws.Cells.Item(29, 1).Value = "04/01/2026 22:45:02"
$ws.Cells.Item(29, 3).Value = "VEJA"
$ws.Cells.Item(29, 4).Value = "Marco Rubio detalha exigências dos EUA para líderes da Venezuela após captura de Maduro"
$ws.Cells.Item(29, 5).Value = "https://veja.abril.com.br/mundo/marco-rubio-detalha-exigencias-dos-eua-para-lideres-da-venezuela-apos-captura-de-maduro/"
$ws.Cells.Item(29, 6).Value = 2
$ws.Cells.Item(29, 7).Value = "senado"
$ws.Cells.Item(29, 8).Value = "Senador americano estabelece condições sobre petróleo, tráfico e grupos armados para q"

# Row 30
$ws.Cells.Item(30, 1).Value = "04/01/2026 22:45:03"
$ws.Cells.Item(30, 3).Value = "VEJA"
$ws.Cells.Item(30, 4).Value = "‘O Agente Secreto’ vence o Critics Choice Awards 2026 de melhor filme internacional"
$ws.Cells.Item(30, 5).Value = "https://veja.abril.com.br/cultura/o-agente-secreto-vence-o-critics-choice-awards-2026-de-melhor-filme-internacional/"
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = "sc"
$ws.Cells.Item(30, 8).Value = "a é a primeira do Brasil a vencer a premiação, considerada &#x27;termômetro&#x27; para o O<b>sc</b>ar"

# Row 31
$ws.Cells.Item(31, 1).Value = "04/01/2026 22:45:03"
$ws.Cells.Item(31, 3).Value = "VEJA"
$ws.Cells.Item(31, 4).Value = "EUA não têm tropas em solo venezuelano, mas mantêm forças no Caribe, diz Pentágono"
$ws.Cells.Item(31, 5).Value = "https://veja.abril.com.br/mundo/estados-unidos-nao-tem-tropas-em-solo-venezuelano-mas-mantem-forcas-no-caribe-diz-pentagono/"
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = "sc"
$ws.Cells.Item(31, 8).Value = "Forças americanas seguem em prontidão no Caribe, enquanto Trump não de<b>sc</b>arta segunda operação militar em Caracas"

